$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the values between row 5 and row 6 for columns A, Q, R, Z, AB
$tmpA = $ws.Range("A5").Value2
$ws.Range("A5").Value2 = $ws.Range("A6").Value2
$ws.Range("A6").Value2 = $tmpA

$tmpQ = $ws.Range("Q5").Value2
$ws.Range("Q5").Value2 = $ws.Range("Q6").Value2
$ws.Range("Q6").Value2 = $tmpQ

$tmpR = $ws.Range("R5").Value2
$ws.Range("R5").Value2 = $ws.Range("R6").Value2
$ws.Range("R6").Value2 = $tmpR

$tmpZ = $ws.Range("Z5").Value2
$ws.Range("Z5").Value2 = $ws.Range("Z6").Value2
$ws.Range("Z6").Value2 = $tmpZ

$tmpAB = $ws.Range("AB5").Value2
$ws.Range("AB5").Value2 = $ws.Range("AB6").Value2
$ws.Range("AB6").Value2 = $tmpAB
